$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DraftOrder")

# Rename header G1: TradeTeam -> TradeWith
$ws.Range("G1").Value = "TradeWith"

# Clear all existing G2:G225 values (previously "nan" placeholder text)
$ws.Range("G2:G225").Value = ""

# Apply trade-up/trade-down changes per row
$ws.Range("F7").Value = "No"
$ws.Range("G7").Value = ""
$ws.Range("F18").Value = "Yes"
$ws.Range("G18").Value = "Chiefs"
$ws.Range("F37").Value = "No"
$ws.Range("G37").Value = ""
$ws.Range("F46").Value = "Yes"
$ws.Range("G46").Value = "Bears"
$ws.Range("F57").Value = "Yes"
$ws.Range("G57").Value = "Panthers"
$ws.Range("F60").Value = "No"
$ws.Range("G60").Value = ""
$ws.Range("F101").Value = "No"
$ws.Range("G101").Value = ""
$ws.Range("F105").Value = "Yes"
$ws.Range("G105").Value = "Raiders"
$ws.Range("F111").Value = "Yes"
$ws.Range("G111").Value = "Bears"
$ws.Range("F117").Value = "No"
$ws.Range("G117").Value = ""
$ws.Range("F119").Value = "Yes"
$ws.Range("G119").Value = "Jets"
$ws.Range("F134").Value = "No"
$ws.Range("G134").Value = ""
$ws.Range("F145").Value = "Yes"
$ws.Range("G145").Value = "Bears"
$ws.Range("F149").Value = "No"
$ws.Range("G149").Value = ""
$ws.Range("F152").Value = "Yes"
$ws.Range("G152").Value = "Seahawks"
$ws.Range("F157").Value = "No"
$ws.Range("G157").Value = ""
$ws.Range("F159").Value = "Yes"
$ws.Range("G159").Value = "Patriots"
$ws.Range("F165").Value = "No"
$ws.Range("G165").Value = ""
$ws.Range("F191").Value = "No"
$ws.Range("G191").Value = ""
$ws.Range("F207").Value = "No"
$ws.Range("G207").Value = ""
$ws.Range("F214").Value = "Yes"
$ws.Range("G214").Value = "Titans"
$ws.Range("F221").Value = "No"
$ws.Range("G221").Value = ""
